$d = $word.ActiveDocument

# The second table in the document ("APELLIDOS Y NOMBRE" / D.N.I / ... roster)
# currently has a header row plus a single student row (Santos Matín-Nieto
# Álvaro). A new student row must be appended below it:
#   Jiménez Coello Daniel | 11a | (empty) | fgbhfxd | 400 | 2022-03-04 | 2022-03-05

$t = $d.Tables.Item(2)

$t.Rows.Add() | Out-Null
$r = $t.Rows.Count

$t.Cell($r, 1).Range.Text = "Jiménez Coello Daniel"
$t.Cell($r, 2).Range.Text = "11a"
$t.Cell($r, 3).Range.Text = ""
$t.Cell($r, 4).Range.Text = "fgbhfxd"
$t.Cell($r, 5).Range.Text = "400"
$t.Cell($r, 6).Range.Text = "2022-03-04"
$t.Cell($r, 7).Range.Text = "2022-03-05"
